$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.009.76'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '1.642.76'
$ws.Range('E3').Value = '  +2.49%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.68'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.28'
$ws.Range('E8').Value = '  +5.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.260'
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').Value = '1.878.56'
$ws.Range('E12').Value = '  +2.53%  '
$ws.Range('D13').Value = '1.654.39'
$ws.Range('E13').Value = '  +3.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.572'
$ws.Range('E14').Value = '  +4.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.38'
$ws.Range('E15').Value = '  +20.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.92'
$ws.Range('E16').Value = '  +4.58%  '
$ws.Range('D17').Value = '30.069.07'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.89'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.01'
$ws.Range('E19').Value = '  +1.72%  '
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.01'
$ws.Range('E22').Value = '  +5.90%  '
$ws.Range('E23').Value = '  +4.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.15'
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.48'
$ws.Range('E25').Value = '  +1.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.73'
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('E27').Value = '  +2.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.69'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0493'
$ws.Range('E30').Value = '  +2.31%  '
$ws.Range('E31').Value = '  +5.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.42'
$ws.Range('E32').Value = '  +5.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.20'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').Value = '1.437.12'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('E35').Value = '  +7.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.05'
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('E39').Value = '  +1.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '77.36'
$ws.Range('E40').Value = '  +16.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.560'
$ws.Range('E41').Value = '  +2.73%  '
$ws.Range('E42').Value = '  +2.46%  '
$ws.Range('E43').Value = '  +3.31%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0497'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '54.83'
$ws.Range('E45').Value = '  -5.68%  '
$ws.Range('E46').Value = '  +6.37%  '
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.41'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('D49').Value = '1.784.70'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('E50').Value = '  +9.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '90.18'
$ws.Range('E51').Value = '  +3.88%  '
